$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-12 22:18:10"
$wsZhCn.Range("E5").Value = "2016-03-12 22:18:10"
$wsZhCn.Range("H4").Value = "2016-03-12 22:18:27"
$wsZhCn.Range("H5").Value = "2016-03-12 22:18:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-12 22:18:13"
$wsDeDe.Range("E5").Value = "2016-03-12 22:18:13"
$wsDeDe.Range("H4").Value = "2016-03-12 22:18:33"
$wsDeDe.Range("H5").Value = "2016-03-12 22:18:33"
